$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Remove the duplicated "Contact | No display for ContactDetail" row (row 11)
$ws.Rows.Item(11).Delete()

# Publisher value
$ws.Range("B9").Value = "Alvearie Team"

# Replace old "Contact" row with new "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive value -> true (now at row 14 after the row deletion above).
# A plain $ws.Range("B14").Value = "true" gets auto-coerced by Excel into the
# Boolean TRUE, not the literal text "true" the source workbook stores. Build
# the text value via a formula in a scratch cell and paste-special just the
# value back in, which keeps it as literal text.
$ws.Range("Z1").Formula = "=""true"""
$ws.Range("Z1").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$excel.CutCopyMode = $false
